$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4670743333333333
$ws.Range("H2").Value = 1.401223
$ws.Range("I2").Value = 0.1980255150811789
$ws.Range("J2").Value = 0.1980255150811789
$ws.Range("M2").Value = 0.5001966666666666
$ws.Range("N2").Value = 1.50059
$ws.Range("O2").Value = 0.03894027965151046
$ws.Range("P2").Value = 0.03894027965151046
$ws.Range("Q2").Value = 0.2336290246188888
$ws.Range("R2").Value = 2.10266122157
$ws.Range("S2").Value = 0.007711168935395511
$ws.Range("T2").Value = 0.007711168935395511
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4670743333333333
$ws.Range("H3").Value = 1.401223
$ws.Range("I3").Value = 0.1980255150811789
$ws.Range("J3").Value = 0.1980255150811789
$ws.Range("O3").Value = 0.7732779360092192
$ws.Range("P3").Value = 0.7732779360092191
$ws.Range("Q3").Value = 4.639416346413777
$ws.Range("R3").Value = 41.75474711772399
$ws.Range("S3").Value = 0.1531287615791366
$ws.Range("T3").Value = 0.1531287615791365
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4670743333333333
$ws.Range("H4").Value = 1.401223
$ws.Range("I4").Value = 0.1980255150811789
$ws.Range("J4").Value = 0.1980255150811789
$ws.Range("M4").Value = 2.334238666666666
$ws.Range("N4").Value = 7.002715999999999
$ws.Range("O4").Value = 0.1817203362411497
$ws.Range("P4").Value = 0.1817203362411496
$ws.Range("Q4").Value = 1.090262969074222
$ws.Range("R4").Value = 9.812366721667999
$ws.Range("S4").Value = 0.03598526318487869
$ws.Range("T4").Value = 0.03598526318487869
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4670743333333333
$ws.Range("H5").Value = 1.401223
$ws.Range("I5").Value = 0.1980255150811789
$ws.Range("J5").Value = 0.1980255150811789
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.07786066666666666
$ws.Range("N5").Value = 0.233582
$ws.Range("O5").Value = 0.006061448098120818
$ws.Range("P5").Value = 0.006061448098120817
$ws.Range("Q5").Value = 0.03636671897622222
$ws.Range("R5").Value = 0.3273004707859999
$ws.Range("S5").Value = 0.001200321381768207
$ws.Range("T5").Value = 0.001200321381768207
$ws.Range("I6").Value = 0.5220544965412525
$ws.Range("J6").Value = 0.5220544965412525
$ws.Range("M6").Value = 0.5001966666666666
$ws.Range("N6").Value = 1.50059
$ws.Range("O6").Value = 0.03894027965151046
$ws.Range("P6").Value = 0.03894027965151046
$ws.Range("Q6").Value = 0.6159159983744443
$ws.Range("R6").Value = 5.543243985369999
$ws.Range("S6").Value = 0.02032894808864487
$ws.Range("T6").Value = 0.02032894808864487
$ws.Range("I7").Value = 0.5220544965412525
$ws.Range("J7").Value = 0.5220544965412525
$ws.Range("O7").Value = 0.7732779360092192
$ws.Range("P7").Value = 0.7732779360092191
$ws.Range("S7").Value = 0.4036932235697518
$ws.Range("T7").Value = 0.4036932235697517
$ws.Range("I8").Value = 0.5220544965412525
$ws.Range("J8").Value = 0.5220544965412525
$ws.Range("M8").Value = 2.334238666666666
$ws.Range("N8").Value = 7.002715999999999
$ws.Range("O8").Value = 0.1817203362411497
$ws.Range("P8").Value = 0.1817203362411496
$ws.Range("Q8").Value = 2.874259335643111
$ws.Range("R8").Value = 25.868334020788
$ws.Range("S8").Value = 0.0948679186476805
$ws.Range("T8").Value = 0.09486791864768049
$ws.Range("I9").Value = 0.5220544965412525
$ws.Range("J9").Value = 0.5220544965412525
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.07786066666666666
$ws.Range("N9").Value = 0.233582
$ws.Range("O9").Value = 0.006061448098120818
$ws.Range("P9").Value = 0.006061448098120817
$ws.Range("Q9").Value = 0.0958735502251111
$ws.Range("R9").Value = 0.8628619520259999
$ws.Range("S9").Value = 0.003164406235175396
$ws.Range("T9").Value = 0.003164406235175395
$ws.Range("G10").Value = 0.6061993333333334
$ws.Range("H10").Value = 1.818598
$ws.Range("I10").Value = 0.2570103443032279
$ws.Range("J10").Value = 0.2570103443032279
$ws.Range("M10").Value = 0.5001966666666666
$ws.Range("N10").Value = 1.50059
$ws.Range("O10").Value = 0.03894027965151046
$ws.Range("P10").Value = 0.03894027965151046
$ws.Range("Q10").Value = 0.3032188858688889
$ws.Range("R10").Value = 2.72896997282
$ws.Range("S10").Value = 0.01000805468049868
$ws.Range("T10").Value = 0.01000805468049868
$ws.Range("G11").Value = 0.6061993333333334
$ws.Range("H11").Value = 1.818598
$ws.Range("I11").Value = 0.2570103443032279
$ws.Range("J11").Value = 0.2570103443032279
$ws.Range("O11").Value = 0.7732779360092192
$ws.Range("P11").Value = 0.7732779360092191
$ws.Range("Q11").Value = 6.021335139913778
$ws.Range("R11").Value = 54.192016259224
$ws.Range("S11").Value = 0.1987404285758189
$ws.Range("T11").Value = 0.1987404285758188
$ws.Range("G12").Value = 0.6061993333333334
$ws.Range("H12").Value = 1.818598
$ws.Range("I12").Value = 0.2570103443032279
$ws.Range("J12").Value = 0.2570103443032279
$ws.Range("M12").Value = 2.334238666666666
$ws.Range("N12").Value = 7.002715999999999
$ws.Range("O12").Value = 0.1817203362411497
$ws.Range("P12").Value = 0.1817203362411496
$ws.Range("Q12").Value = 1.415013923574222
$ws.Range("R12").Value = 12.735125312168
$ws.Range("S12").Value = 0.04670400618423622
$ws.Range("T12").Value = 0.04670400618423621
$ws.Range("G13").Value = 0.6061993333333334
$ws.Range("H13").Value = 1.818598
$ws.Range("I13").Value = 0.2570103443032279
$ws.Range("J13").Value = 0.2570103443032279
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.07786066666666666
$ws.Range("N13").Value = 0.233582
$ws.Range("O13").Value = 0.006061448098120818
$ws.Range("P13").Value = 0.006061448098120817
$ws.Range("Q13").Value = 0.04719908422622222
$ws.Range("R13").Value = 0.424791758036
$ws.Range("S13").Value = 0.001557854862674177
$ws.Range("T13").Value = 0.001557854862674177
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.054036
$ws.Range("H14").Value = 0.162108
$ws.Range("I14").Value = 0.0229096440743406
$ws.Range("J14").Value = 0.0229096440743406
$ws.Range("M14").Value = 0.5001966666666666
$ws.Range("N14").Value = 1.50059
$ws.Range("O14").Value = 0.03894027965151046
$ws.Range("P14").Value = 0.03894027965151046
$ws.Range("Q14").Value = 0.02702862708
$ws.Range("R14").Value = 0.24325764372
$ws.Range("S14").Value = 0.0008921079469713925
$ws.Range("T14").Value = 0.0008921079469713925
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.054036
$ws.Range("H15").Value = 0.162108
$ws.Range("I15").Value = 0.0229096440743406
$ws.Range("J15").Value = 0.0229096440743406
$ws.Range("O15").Value = 0.7732779360092192
$ws.Range("P15").Value = 0.7732779360092191
$ws.Range("Q15").Value = 0.536735769456
$ws.Range("R15").Value = 4.830621925103999
$ws.Range("S15").Value = 0.01771552228451194
$ws.Range("T15").Value = 0.01771552228451194
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.054036
$ws.Range("H16").Value = 0.162108
$ws.Range("I16").Value = 0.0229096440743406
$ws.Range("J16").Value = 0.0229096440743406
$ws.Range("M16").Value = 2.334238666666666
$ws.Range("N16").Value = 7.002715999999999
$ws.Range("O16").Value = 0.1817203362411497
$ws.Range("P16").Value = 0.1817203362411496
$ws.Range("Q16").Value = 0.126132920592
$ws.Range("R16").Value = 1.135196285328
$ws.Range("S16").Value = 0.004163148224354236
$ws.Range("T16").Value = 0.004163148224354235
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.054036
$ws.Range("H17").Value = 0.162108
$ws.Range("I17").Value = 0.0229096440743406
$ws.Range("J17").Value = 0.0229096440743406
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.07786066666666666
$ws.Range("N17").Value = 0.233582
$ws.Range("O17").Value = 0.006061448098120818
$ws.Range("P17").Value = 0.006061448098120817
$ws.Range("Q17").Value = 0.004207278984
$ws.Range("R17").Value = 0.037865510856
$ws.Range("S17").Value = 0.0001388656185030367
$ws.Range("T17").Value = 0.0001388656185030367
